$wb = $excel.ActiveWorkbook

$tools = $wb.Worksheets.Item("Tools")
$sources = $wb.Worksheets.Item("Sources")

# Add two new rows of tool/source info to the "Tools" sheet.
# Cell values are written in this specific order so that new shared-string
# table entries come out in the same order as the authored workbook.
$tools.Range("B8").Value = "Get data about people and the communities they live in, includes population, identity, housing, people in or out of work, education and health."
$tools.Range("B9").Value = "Experimental statistics release showing clustering analysis exploring similarities between local authorities in England"
$tools.Range("C8").Value = "<a href='https://www.ons.gov.uk/visualisations/areas/'>ONS</a>"
$tools.Range("A8").Value = "ONS area information"
$tools.Range("C9").Value = "<a href='https://www.ons.gov.uk/peoplepopulationandcommunity/wellbeing/articles/clusteringlocalauthoritiesagainstsubnationalindicatorsengland/2023-02-24'>ONS</a>"
$tools.Range("A9").Value = "ONS clustering local authorities"

$tools.Range("D8").Value = "Publicly available "
$tools.Range("D9").Value = "Publicly available "

# Printer setup for the "Tools" sheet (portrait, A4/letter-size "9" = A4 paper)
$tools.PageSetup.PaperSize = 9
$tools.PageSetup.Orientation = 1

# "Sources" sheet keeps its own cursor position but is no longer the active tab.
$sources.Range("A13").Select()

# "Tools" sheet becomes the active / selected tab, with its cursor moved past
# the newly added rows. Activated last so it ends up as the visible tab.
$tools.Activate()
$tools.Range("A10").Select()
